$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.476.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.11%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.489.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.96%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'569.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.48%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'163.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -1.59%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.90%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.488.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.01%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.98%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "'  -0.53%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.352"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -1.19%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "'  -1.32%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'2.944.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.08%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'69.321.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.04%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.0000174"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.64%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'24.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.11%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.496.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.87%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'11.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.64%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -4.50%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'346.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.73%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'3.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.31%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").Value = "'SuiNetwork"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'1.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.98%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").Value = "'Dai"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.09%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'69.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.04%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'3.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.65%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.616.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.31%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'8.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.42%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'1.01"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.76%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.0₃0865"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.99%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'7.55"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -4.22%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'437.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -5.65%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.47%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -0.08%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  -1.62%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'156.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.61%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -3.05%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'19.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.16%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'18.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.21%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +0.01%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  -2.14%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'4.56"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.08%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'1.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.96%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +40.81%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'2.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -5.07%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "'  -5.89%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'137.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.60%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'3.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.03%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.504"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -4.32%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0723"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.88%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.572"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.74%  "
$ws.Range("E51").Style = "Normal"
